$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data
$ws.Range("D2").Value = "26.443.67"
$ws.Range("E2").Value = "  -3.78%  "
$ws.Range("D3").Value = "1.772.28"
$ws.Range("E3").Value = "  -2.82%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'1.003"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "'306.51"
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("D7").Value = "'0.4290"
$ws.Range("E7").Value = "  +0.90%  "
$ws.Range("D8").Value = "'0.3660"
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("D9").Value = "'0.07287"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").Value = "'0.8476"
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").Value = "'20.33"
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("D12").Value = "1.761.02"
$ws.Range("E12").Value = "  -4.85%  "
$ws.Range("D13").Value = "'5.257"
$ws.Range("E13").Value = "  -2.57%  "
$ws.Range("D14").Value = "'6.436"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").Value = "'0.06822"
$ws.Range("E15").Value = "  -1.58%  "
$ws.Range("D16").Value = "'1.007"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "'79.52"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "'0.000008700"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").Value = "'1.003"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("D21").Value = "26.449.23"
$ws.Range("E21").Value = "  -4.09%  "
$ws.Range("D22").Value = "'5.098"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").Value = "'11.25"
$ws.Range("E23").Value = "  +3.42%  "
$ws.Range("D24").Value = "2.002.34"
$ws.Range("E24").Value = "  -2.30%  "
$ws.Range("D25").Value = "'152.52"
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("D26").Value = "'1.853"
$ws.Range("E26").Value = "  -6.97%  "
$ws.Range("D27").Value = "'18.18"
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").Value = "'5.092"
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("D29").Value = "'114.73"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").Value = "'1.711"
$ws.Range("E30").Value = "  -4.68%  "
$ws.Range("D31").Value = "'0.08949"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("D32").Value = "'0.7259"
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("D33").Value = "'1.115"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("D34").Value = "'4.337"
$ws.Range("E34").Value = "  -4.40%  "
$ws.Range("D35").Value = "'2.758"
$ws.Range("E35").Value = "  -7.19%  "
$ws.Range("D36").Value = "'1.003"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("D38").Value = "'0.05155"
$ws.Range("E38").Value = "  -2.42%  "
$ws.Range("D39").Value = "'0.01894"
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("D40").Value = "'0.4927"
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("D41").Value = "'0.1610"
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("D42").Value = "'2.518"
$ws.Range("E42").Value = "  -9.99%  "
$ws.Range("D43").Value = "'6.214"
$ws.Range("E43").Value = "  -3.69%  "
$ws.Range("D44").Value = "'8.062"
$ws.Range("E44").Value = "  -3.50%  "
$ws.Range("D45").Value = "'104.90"
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "'10.12"
$ws.Range("E47").Value = "  -3.53%  "
$ws.Range("D48").Value = "'0.4502"
$ws.Range("E48").Value = "  -3.72%  "
$ws.Range("D49").Value = "'0.06198"
$ws.Range("E49").Value = "  -4.25%  "
$ws.Range("D50").Value = "'1.581"
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("D51").Value = "'1.743"
$ws.Range("E51").Value = "  +2.58%  "
